# Add "marker_1" column (J) with marker values per strain/row to the bioSample sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column J
$ws.Cells.Item(1, 10).Value = "marker_1"

# Marker values for each data row (J4..J25), written explicitly per row.
# Rows 2, 3, 15, 22, 23, 24 are intentionally left blank (no marker).
$ws.Cells.Item(4, 10).Value = "G418"
$ws.Cells.Item(5, 10).Value = "NAT"
$ws.Cells.Item(6, 10).Value = "NAT"
$ws.Cells.Item(7, 10).Value = "NAT"
$ws.Cells.Item(8, 10).Value = "NAT"
$ws.Cells.Item(9, 10).Value = "NAT"
$ws.Cells.Item(10, 10).Value = "NAT"
$ws.Cells.Item(11, 10).Value = "G418"
$ws.Cells.Item(12, 10).Value = "G418"
$ws.Cells.Item(13, 10).Value = "NAT"
$ws.Cells.Item(14, 10).Value = "G418"
$ws.Cells.Item(16, 10).Value = "NAT"
$ws.Cells.Item(17, 10).Value = "NAT"
$ws.Cells.Item(18, 10).Value = "NAT"
$ws.Cells.Item(19, 10).Value = "NAT"
$ws.Cells.Item(20, 10).Value = "NAT"
$ws.Cells.Item(21, 10).Value = "NAT"
$ws.Cells.Item(25, 10).Value = "G418"

# Update the sheet view selection / scroll position to match the authored state
$ws.Range("J26").Select()
$excel.ActiveWindow.ScrollColumn = 2
